$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5909176666666667
$ws.Range("H2").Value = 1.772753
$ws.Range("I2").Value = 0.9937758428931484
$ws.Range("J2").Value = 0.9937758428931482
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 1.005700307457444
$ws.Range("R2").Value = 9.051302767116999
$ws.Range("S2").Value = 0.02092517159477661
$ws.Range("T2").Value = 0.02092517159477661

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5909176666666667
$ws.Range("H3").Value = 1.772753
$ws.Range("I3").Value = 0.9937758428931484
$ws.Range("J3").Value = 0.9937758428931482
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 62.503947
$ws.Range("N3").Value = 187.511841
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 36.934686518697
$ws.Range("R3").Value = 332.412178668273
$ws.Range("S3").Value = 0.7684840577974273
$ws.Range("T3").Value = 0.7684840577974273

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5909176666666667
$ws.Range("H4").Value = 1.772753
$ws.Range("I4").Value = 0.9937758428931484
$ws.Range("J4").Value = 0.9937758428931482
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 0.256380857119
$ws.Range("R4").Value = 2.307427714071
$ws.Range("S4").Value = 0.005334405676372918
$ws.Range("T4").Value = 0.005334405676372918

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.5909176666666667
$ws.Range("H5").Value = 1.772753
$ws.Range("I5").Value = 0.9937758428931484
$ws.Range("J5").Value = 0.9937758428931482
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 9.438363884384
$ws.Range("R5").Value = 84.945274959456
$ws.Range("S5").Value = 0.1963799577172095
$ws.Range("T5").Value = 0.1963799577172095

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.5909176666666667
$ws.Range("H6").Value = 1.772753
$ws.Range("I6").Value = 0.9937758428931484
$ws.Range("J6").Value = 0.9937758428931482
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 0.1274717741905556
$ws.Range("R6").Value = 1.147245967715
$ws.Range("S6").Value = 0.002652250107361988
$ws.Range("T6").Value = 0.002652250107361988

$ws.Range("I7").Value = 0.006224157106851674
$ws.Range("J7").Value = 0.006224157106851673
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 0.006298841696333332
$ws.Range("R7").Value = 0.056689575267
$ws.Range("S7").Value = 0.0001310572765731067
$ws.Range("T7").Value = 0.0001310572765731067

$ws.Range("I8").Value = 0.006224157106851674
$ws.Range("J8").Value = 0.006224157106851673
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("S8").Value = 0.004813123144467862
$ws.Range("T8").Value = 0.004813123144467862

$ws.Range("I9").Value = 0.006224157106851674
$ws.Range("J9").Value = 0.006224157106851673
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 0.001605749169
$ws.Range("R9").Value = 0.014451742521
$ws.Range("S9").Value = 0.00003341012889261421
$ws.Range("T9").Value = 0.00003341012889261421

$ws.Range("I10").Value = 0.006224157106851674
$ws.Range("J10").Value = 0.006224157106851673
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 0.059113793184
$ws.Range("R10").Value = 0.5320241386560001
$ws.Range("S10").Value = 0.001229955143516427
$ws.Range("T10").Value = 0.001229955143516427

$ws.Range("I11").Value = 0.006224157106851674
$ws.Range("J11").Value = 0.006224157106851673
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 0.0007983735516666667
$ws.Range("R11").Value = 0.007185361965
$ws.Range("S11").Value = 0.00001661141340166405
$ws.Range("T11").Value = 0.00001661141340166405
